# Update "想去人数" (column F) counts scraped for 北京-漫展信息.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(6, 6).Value = 1601
$ws.Cells.Item(9, 6).Value = 732
$ws.Cells.Item(10, 6).Value = 2685
$ws.Cells.Item(11, 6).Value = 2685
$ws.Cells.Item(13, 6).Value = 1777
$ws.Cells.Item(14, 6).Value = 610
$ws.Cells.Item(16, 6).Value = 698
$ws.Cells.Item(17, 6).Value = 5097
$ws.Cells.Item(18, 6).Value = 231
$ws.Cells.Item(19, 6).Value = 74
$ws.Cells.Item(21, 6).Value = 3389
$ws.Cells.Item(22, 6).Value = 865
$ws.Cells.Item(24, 6).Value = 71
$ws.Cells.Item(25, 6).Value = 41
$ws.Cells.Item(26, 6).Value = 2431
$ws.Cells.Item(28, 6).Value = 371
$ws.Cells.Item(31, 6).Value = 483
$ws.Cells.Item(32, 6).Value = 1304
$ws.Cells.Item(33, 6).Value = 804
$ws.Cells.Item(34, 6).Value = 7
$ws.Cells.Item(36, 6).Value = 23
$ws.Cells.Item(37, 6).Value = 56
$ws.Cells.Item(38, 6).Value = 1453
$ws.Cells.Item(39, 6).Value = 19
$ws.Cells.Item(40, 6).Value = 1402
$ws.Cells.Item(41, 6).Value = 91

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(7, 6).Value = 95
$ws.Cells.Item(17, 6).Value = 333
$ws.Cells.Item(18, 6).Value = 260

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 886
$ws.Cells.Item(4, 6).Value = 250
$ws.Cells.Item(6, 6).Value = 39
$ws.Cells.Item(7, 6).Value = 62
$ws.Cells.Item(8, 6).Value = 7

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 886
$ws.Cells.Item(5, 6).Value = 250
$ws.Cells.Item(8, 6).Value = 39
$ws.Cells.Item(9, 6).Value = 62
$ws.Cells.Item(13, 6).Value = 1601
$ws.Cells.Item(17, 6).Value = 2685
$ws.Cells.Item(18, 6).Value = 7
$ws.Cells.Item(20, 6).Value = 1777
$ws.Cells.Item(22, 6).Value = 610
$ws.Cells.Item(24, 6).Value = 698
$ws.Cells.Item(25, 6).Value = 5097
$ws.Cells.Item(26, 6).Value = 231
$ws.Cells.Item(27, 6).Value = 74
$ws.Cells.Item(29, 6).Value = 3389
$ws.Cells.Item(30, 6).Value = 865
$ws.Cells.Item(32, 6).Value = 71
$ws.Cells.Item(34, 6).Value = 41
$ws.Cells.Item(35, 6).Value = 2431
$ws.Cells.Item(37, 6).Value = 371
$ws.Cells.Item(39, 6).Value = 483
$ws.Cells.Item(40, 6).Value = 1304
$ws.Cells.Item(42, 6).Value = 260
$ws.Cells.Item(44, 6).Value = 804
$ws.Cells.Item(45, 6).Value = 7
$ws.Cells.Item(47, 6).Value = 23
$ws.Cells.Item(48, 6).Value = 56
$ws.Cells.Item(49, 6).Value = 1402
$ws.Cells.Item(50, 6).Value = 91
